$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value that gets refreshed on each
# automatic run. Update it from 2023-10-05 (45204) to 2023-10-08 (45207)
# for every data row (rows 2 through 97).
$ws.Range("C2:C97").Value = 45207
